# Packages zugeordnet, bom ergaenzt.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Neue Zeilen fuer zugeordnete Packages einfuegen (Taster, Knopfzellenhalter) ---
$ws.Range("A33").Value = 1
$ws.Range("B33").Value = "Taster"
$ws.Range("D33").Value = "611-PTS636SM25JSMTRL "

$ws.Range("A34").Value = 1
$ws.Range("B34").Value = "Knopfzellenhalter"
$ws.Range("D34").Value = "534-500 "

# --- Bestehenden Power-Supply Block (Zeilen 37:38) ans Ende der BOM verschieben (Zeilen 48:49) ---
$ws.Rows("37:38").Cut()
$ws.Rows("48").Select()
$ws.Paste()
$excel.CutCopyMode = $false

# --- Sicht aktualisieren (Zoom + Scrollposition + Auswahl) ---
$win = $excel.ActiveWindow
$win.Zoom = 55
$ws.Range("B36").Select()
